$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41.. down by one
# (Excel will extend the used range / dimension automatically, e.g. to A1:R136)
$ws.Rows("41").Insert()

# Populate the newly inserted row 41 with the new weekly data point.
# (Columns A,B,C,E,F,G,H,I,J,N,Q,R mirror the same record that used to sit
# at row 41 and now lives at row 42; D,K,L,M,O,P hold the new values.)
$ws.Range("A41").Value2 = 7
$ws.Range("B41").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C41").Value = "Ñuble"
$ws.Range("D41").Value2 = 44414
$ws.Range("E41").Value2 = 16
$ws.Range("F41").Value2 = 100112023
$ws.Range("G41").Value = "Brócoli"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value2 = 300
$ws.Range("K41").Value2 = 650
$ws.Range("L41").Value2 = 700
$ws.Range("M41").Value2 = 675
$ws.Range("N41").Value = "$/unidad"
$ws.Range("O41").Value = "Provincia de Diguillín"
$ws.Range("P41").Value2 = 675
$ws.Range("Q41").Value2 = 1
$ws.Range("R41").Value = "Hortaliza"
